$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.248.38'
$ws.Range("E2").Value = '  +0.29%  '
$ws.Range("D3").Value = '1.650.88'
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("E6").Value = '  +1.46%  '
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("E9").Value = '  +0.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.01'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.38%  '
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("D12").Value = '1.881.63'
$ws.Range("E12").Value = '  -0.30%  '
$ws.Range("D13").Value = '1.646.02'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.85%  '
$ws.Range("D17").Value = '27.236.13'
$ws.Range("E17").Value = '  +0.16%  '
$ws.Range("E18").Value = '  +0.86%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '220.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("E20").Value = '  -0.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.83'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.91%  '
$ws.Range("E22").Value = '  +6.23%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("E24").Value = '  -0.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.85'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("E27").Value = '  +1.23%  '
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.16%  '
$ws.Range("E30").Value = '  -0.95%  '
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.04'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.58'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.88%  '
$ws.Range("D35").Value = '1.259.89'
$ws.Range("E35").Value = '  +0.15%  '
$ws.Range("E36").Value = '  +0.28%  '
$ws.Range("E37").Value = '  +0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.548'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.14%  '
$ws.Range("E39").Value = '  +0.81%  '
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.810'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("E42").Value = '  +1.34%  '
$ws.Range("E43").Value = '  +4.06%  '
$ws.Range("D44").Value = '1.792.44'
$ws.Range("E44").Value = '  -0.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.72'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("E47").Value = '  -0.47%  '
$ws.Range("D48").Value = '0.0₆0107'
$ws.Range("E48").Value = '  +17.41%  '
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.69'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.21%  '
$ws.Range("E51").Value = '  -1.24%  '
